# ---------------------------------------------------------------------------
# Obsidian vault auto-backup edit:
#   - sheet1 ("Feuil1"): F25 now returns a positive PV (negated PV formula),
#     selection moves from F25 to F26, sheet is no longer the active tab.
#   - a new sheet "20M ma gueule" is appended after "Feuil1" and becomes the
#     active tab; it reproduces the rente/VA table with B2 = 1 000 000 and a
#     much larger rate sweep (1%..20%), plus a sensitivity grid on row 25/26
#     and a small "Pourcentage idéal" goal-seek style block at F28:F30.
# ---------------------------------------------------------------------------

function ColLetter([int]$n) {
    $letters = ""
    while ($n -gt 0) {
        $rem = ($n - 1) % 26
        $letters = [char](65 + $rem) + $letters
        $n = [int](($n - $rem - 1) / 26)
    }
    return $letters
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) sheet1 ("Feuil1") edits
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F25").Formula = "=-1*PV(B3,A25,B2)"
$ws1.Range("F26").Select()

# ---------------------------------------------------------------------------
# 2) Add the new sheet right after Feuil1 -> becomes sheetId 2, active tab
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "20M ma gueule"

# ---- Row 2 : Rente ----------------------------------------------------
$ws2.Range("A2").Value = "Rente"
$ws2.Range("B2").Value = 1000000

# ---- Row 3 : Taux d'intérêt sweep (B..U = 1% .. 20%) -------------------
$ws2.Range("A3").Value = "Taux d'intérêt"
for ($i = 0; $i -lt 20; $i++) {
    $col = ColLetter (2 + $i)
    $cell = $ws2.Range("$col`3")
    $cell.Value = [math]::Round(0.01 * ($i + 1), 2)
    $cell.NumberFormat = "0%"
}

# ---- Row 5 : table header ----------------------------------------------
$ws2.Range("A5").Value = "Periode"
$ws2.Range("B5").Value = "Rente"
$ws2.Range("C5").Value = "VA"
$ws2.Range("D5").Value = "VA Cumulée"

# ---- Rows 6..25 : period table -----------------------------------------
for ($r = 6; $r -le 25; $r++) {
    $a = $r - 5
    $ws2.Range("A$r").Value = $a

    $bCell = $ws2.Range("B$r")
    $bCell.Formula = "=`$B`$2"
    $bCell.Interior.Color = 65535

    $cCell = $ws2.Range("C$r")
    $cCell.Formula = "=`$B$r/((1 + `$B`$3)^A$r)"
    $cCell.Interior.Color = 65535

    $dCell = $ws2.Range("D$r")
    $dCell.Formula = "=SUM(`$C`$6:C$r)"
    $dCell.Interior.Color = 65535
}

$ws2.Range("F24").Value = "Formule"

# ---- Row 25 F:Y : sensitivity PV grid (uses row-3 rates, B..U) --------
for ($i = 0; $i -lt 20; $i++) {
    $col = ColLetter (6 + $i)
    $rateCol = ColLetter (2 + $i)
    $cell = $ws2.Range("$col`25")
    $cell.Formula = "=-1*PV($rateCol`$3,`$A`$25,`$B`$2)"
    $cell.Interior.Color = 65535
    $cell.NumberFormat = "#,##0.00 ""€""_);[Red]\(#,##0.00 ""€""\)"
}

# ---- Row 26 F:Y : same rate sweep, shifted to columns F..Y -------------
for ($i = 0; $i -lt 20; $i++) {
    $col = ColLetter (6 + $i)
    $cell = $ws2.Range("$col`26")
    $cell.Value = [math]::Round(0.01 * ($i + 1), 2)
    $cell.NumberFormat = "0%"
}

# ---- "Pourcentage idéal" block -----------------------------------------
$ws2.Range("F28").Value = "Pourcentage idéal"
$ws2.Range("F29").Formula = "=-1*PV(F30,`$A`$25,`$B`$2)"
$ws2.Range("F30").Value = 0.077546895300093871
$ws2.Range("F30").NumberFormat = "0%"

$ws2.Range("F29").Select()
$excel.ActiveWindow.Zoom = 125
